# Update the build timestamp embedded in version / citation strings
# from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet -------------------------------------------------------

$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")"

$wsAbout.Range("A6").Value = "Recommended Citation:  " + [char]34 + "Global Energy Monitor, Coal mine boundaries and methane sources for Leer Coal Mine, United States, M1036, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet -------------------------------
# Column S ("build_version") for data rows 2 through 11.

for ($row = 2; $row -le 11; $row++) {
    $cell = $wsData.Range("S" + $row)
    $cell.Value = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on " + $newStamp + ")"
}
